# Update "想去人数" (interest count) figures in column F for the
# "展览" sheet and the corresponding rows in the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 14178
    3  = 561
    5  = 1223
    6  = 1057
    7  = 13973
    8  = 15151
    10 = 29
    12 = 180
    18 = 27
    19 = 67
    20 = 24
    21 = 1178
    24 = 5884
    25 = 953
    26 = 1073
    27 = 5484
    30 = 77
    31 = 382
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 14178
    4  = 561
    6  = 1223
    7  = 1057
    8  = 13973
    9  = 15151
    11 = 29
    13 = 180
    19 = 27
    20 = 67
    21 = 24
    22 = 1178
    26 = 5884
    27 = 953
    28 = 1073
    29 = 5484
    32 = 77
    33 = 382
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
